$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (which currently holds "nama"),
# shifting nama/gender/Kelas one column to the right.
$ws.Range("C1").EntireColumn.Insert()

# Set the new header cell's value and copy the header style (bold) from A1.
$ws.Range("C1").Value = "nisn"
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)  # xlPasteFormats

# Update the active cell selection to A2, matching the target state.
$ws.Range("A2").Select()
